# Add a new worksheet "USERACCOUNTMANAGEMENTDATA" with the same header
# as the DATA sheet plus the two "change password" rows, mirroring the
# rows already present in DATA (rows 10-11).

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("DATA")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Create the new sheet after the last existing sheet.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "USERACCOUNTMANAGEMENTDATA"

# Header row (copy of DATA's header row).
$headers = @("testname", "execute", "username", "password", "fname", "browser", "currentpassword", "newpassword", "confirmpassword")
for ($c = 1; $c -le $headers.Length; $c++) {
    $newSheet.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Data rows - identical to DATA sheet rows 10 and 11 (the
# "verifyThatUserCanChangePasswordWithValidCredentials" rows), one per browser.
$row2 = @("verifyThatUserCanChangePasswordWithValidCredentials", "yes", "Admin", "admin123", "amuthan", "chrome", "admin123", "admin111", "admin111")
$row3 = @("verifyThatUserCanChangePasswordWithValidCredentials", "yes", "Admin", "admin123", "amuthan", "firefox", "admin123", "admin111", "admin111")

for ($c = 1; $c -le $row2.Length; $c++) {
    $newSheet.Cells.Item(2, $c).Value = $row2[$c - 1]
}
for ($c = 1; $c -le $row3.Length; $c++) {
    $newSheet.Cells.Item(3, $c).Value = $row3[$c - 1]
}

# Column widths matching the DATA sheet layout.
$newSheet.Columns.Item(1).ColumnWidth = 45.578125
$newSheet.Columns.Item(2).ColumnWidth = 6.83984375
$newSheet.Columns.Item(3).ColumnWidth = 8.5234375
$newSheet.Columns.Item(4).ColumnWidth = 8.62890625
$newSheet.Columns.Item(5).ColumnWidth = 7.7890625
$newSheet.Columns.Item(6).ColumnWidth = 7.15625
$newSheet.Columns.Item(7).ColumnWidth = 14.05078125
$newSheet.Columns.Item(8).ColumnWidth = 11.5234375
$newSheet.Columns.Item(9).ColumnWidth = 14.41796875

# Selection state matching the target workbook.
$newSheet.Range("A2:I3").Select()
$newSheet.Range("A2").Activate()

$dataSheet.Activate()
$dataSheet.Range("H6").Select()
